$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 146
$ws.Range("I11").Value = 146
$ws.Range("K11").Value = 146
$ws.Range("M11").Value = -6
$ws.Range("H19").Value = 810
$ws.Range("I19").Value = 585
$ws.Range("J19").Value = 960
$ws.Range("K19").Value = 585
$ws.Range("L19").Value = 960
$ws.Range("M19").Value = -410
$ws.Range("N19").Value = -1310
$ws.Range("H40").Value = 4284.375
$ws.Range("J40").Value = 5812.625
$ws.Range("L40").Value = 5812.625
$ws.Range("N40").Value = -6162.625
$ws.Range("H111").Value = 1803
$ws.Range("I111").Value = 1074.5
$ws.Range("J111").Value = 2531.5
$ws.Range("K111").Value = 3223.5
$ws.Range("L111").Value = 7594.5
$ws.Range("M111").Value = -156.5
$ws.Range("N111").Value = -13728.5
$ws.Range("H137").Value = 1077.3636
$ws.Range("I137").Value = 699.5
$ws.Range("J137").Value = 1161.3334
$ws.Range("K137").Value = 2098.5
$ws.Range("L137").Value = 3484.0002
$ws.Range("M137").Value = 451.5
$ws.Range("N137").Value = -8584.0002
$ws.Range("H138").Value = 3647.348
$ws.Range("I138").Value = 2174.5
$ws.Range("J138").Value = 4432.8667
$ws.Range("K138").Value = 6523.5
$ws.Range("L138").Value = 13298.6001
$ws.Range("M138").Value = -1383.5
$ws.Range("N138").Value = -23578.6001
$ws.Range("H141").Value = 1460
$ws.Range("I141").Value = 1460
$ws.Range("K141").Value = 4380
$ws.Range("M141").Value = 800

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 463.33334
$ws.Range("I4").Value = 445
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 445
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -329
$ws.Range("N4").Value = -732
$ws.Range("H5").Value = 510.75
$ws.Range("I5").Value = 671
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 671
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = -559
$ws.Range("N5").Value = -254
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 5000
$ws.Range("K22").Value = 5000
$ws.Range("M22").Value = -4701
$ws.Range("H61").Value = 13758.25
$ws.Range("I61").Value = 10010.4
$ws.Range("J61").Value = 20004.666
$ws.Range("K61").Value = 10010.4
$ws.Range("L61").Value = 20004.666
$ws.Range("M61").Value = -9798.4
$ws.Range("N61").Value = -20428.666
$ws.Range("H74").Value = 2463.6365
$ws.Range("I74").Value = 2463.6365
$ws.Range("K74").Value = 2463.6365
$ws.Range("M74").Value = -1589.6365
$ws.Range("H77").Value = 2463.6365
$ws.Range("I77").Value = 2463.6365
$ws.Range("K77").Value = 12318.1825
$ws.Range("M77").Value = -7950.182500000001
$ws.Range("H132").Value = 3499.1667
$ws.Range("I132").Value = 3499.1667
$ws.Range("K132").Value = 10497.5001
$ws.Range("M132").Value = -7967.500100000001
$ws.Range("H136").Value = 13758.25
$ws.Range("I136").Value = 10010.4
$ws.Range("J136").Value = 20004.666
$ws.Range("K136").Value = 30031.2
$ws.Range("L136").Value = 60013.99800000001
$ws.Range("M136").Value = -27481.2
$ws.Range("N136").Value = -65113.99800000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 510.75
$ws.Range("I4").Value = 671
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 671
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = -556
$ws.Range("N4").Value = -260
$ws.Range("H22").Value = 785.2857
$ws.Range("I22").Value = 749.6667
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 749.6667
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -576.6667
$ws.Range("N22").Value = -1345
$ws.Range("H24").Value = 804
$ws.Range("I24").Value = 608
$ws.Range("J24").Value = 1000
$ws.Range("K24").Value = 608
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = -373
$ws.Range("N24").Value = -1470
$ws.Range("H94").Value = 1082.5454
$ws.Range("I94").Value = 1110.8
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 1110.8
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -659.8
$ws.Range("N94").Value = -1702
$ws.Range("H134").Value = 941.7143
$ws.Range("I134").Value = 941.7143
$ws.Range("K134").Value = 2825.1429
$ws.Range("M134").Value = -290.1428999999998

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1882.4667
$ws.Range("I7").Value = 1304.1111
$ws.Range("J7").Value = 2750
$ws.Range("K7").Value = 1304.1111
$ws.Range("L7").Value = 2750
$ws.Range("M7").Value = -1191.1111
$ws.Range("N7").Value = -2976
$ws.Range("H22").Value = 4795.0713
$ws.Range("J22").Value = 700
$ws.Range("L22").Value = 700
$ws.Range("N22").Value = -1400
$ws.Range("H35").Value = 2868.5
$ws.Range("I35").Value = 2691.3333
$ws.Range("K35").Value = 2691.3333
$ws.Range("M35").Value = -2397.3333

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 29.53125
$ws.Range("I2").Value = 15.666667
$ws.Range("J2").Value = 41.764706
$ws.Range("K2").Value = 94.00000199999999
$ws.Range("L2").Value = 250.588236
$ws.Range("M2").Value = 18.99999800000001
$ws.Range("N2").Value = -476.588236
$ws.Range("H10").Value = 112.75
$ws.Range("I10").Value = 112.75
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 338.25
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -199.25
$ws.Range("N10").ClearContents()
$ws.Range("H12").Value = 371.33334
$ws.Range("J12").Value = 423.6154
$ws.Range("L12").Value = 1270.8462
$ws.Range("N12").Value = -1616.8462
$ws.Range("H17").Value = 299.66666
$ws.Range("J17").Value = 699
$ws.Range("L17").Value = 2097
$ws.Range("N17").Value = -2435
$ws.Range("H38").Value = 156.125
$ws.Range("I38").Value = 53
$ws.Range("J38").Value = 218
$ws.Range("K38").Value = 159
$ws.Range("L38").Value = 654
$ws.Range("M38").Value = 188
$ws.Range("N38").Value = -1348
$ws.Range("H51").Value = 1833.3334
$ws.Range("I51").Value = 1500
$ws.Range("K51").Value = 4500
$ws.Range("M51").Value = -4040
$ws.Range("H108").Value = 321.25
$ws.Range("I108").Value = 321.25
$ws.Range("K108").Value = 963.75
$ws.Range("M108").Value = 1916.25
$ws.Range("H132").Value = 1750
$ws.Range("I132").Value = 500
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -32060

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 137.86667
$ws.Range("I2").Value = 64.875
$ws.Range("J2").Value = 221.28572
$ws.Range("K2").Value = 64.875
$ws.Range("L2").Value = 221.28572
$ws.Range("M2").Value = 48.125
$ws.Range("N2").Value = -447.28572
$ws.Range("H14").Value = 1030000
$ws.Range("I14").Value = 2500000
$ws.Range("J14").Value = 50000
$ws.Range("K14").Value = 2500000
$ws.Range("L14").Value = 50000
$ws.Range("M14").Value = -2499832
$ws.Range("N14").Value = -50336

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1857.1428
$ws.Range("I22").Value = 600.3333
$ws.Range("J22").Value = 2799.75
$ws.Range("K22").Value = 600.3333
$ws.Range("L22").Value = 2799.75
$ws.Range("M22").Value = -305.3333
$ws.Range("N22").Value = -3389.75
$ws.Range("H27").Value = 1857.1428
$ws.Range("I27").Value = 600.3333
$ws.Range("J27").Value = 2799.75
$ws.Range("K27").Value = 600.3333
$ws.Range("L27").Value = 2799.75
$ws.Range("M27").Value = -493.3333
$ws.Range("N27").Value = -3013.75
$ws.Range("H132").Value = 4644.7
$ws.Range("I132").Value = 4560
$ws.Range("K132").Value = 13680
$ws.Range("M132").Value = -11150

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H136").Value = 1014.4667
$ws.Range("J136").Value = 1166.3334
$ws.Range("L136").Value = 3499.0002
$ws.Range("N136").Value = -8599.0002
